$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 147.96
$ws.Range("B3").Value = 71.38
$ws.Range("B4").Value = 202.83
$ws.Range("B5").Value = 0.1475
$ws.Range("B6").Value = 0.285
$ws.Range("B7").Value = 0.425
$ws.Range("B8").Value = 270.41
$ws.Range("B9").Value = 121.18
$ws.Range("B10").Value = 2344.76
$ws.Range("B11").Value = 196.64

$ws.Range("B2:B11").Select()
